$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Un-hide rows that were previously filtered/hidden
$ws.Rows.Item(2).Hidden = $false
$ws.Rows.Item(3).Hidden = $false
$ws.Rows.Item(7).Hidden = $false

# Row 8 got a slightly smaller, explicit custom height
$ws.Rows.Item(8).RowHeight = 71.25

# New row 10: common_pitfalls.ipynb / code quality notebook
$ws.Range("B10").Value = "code quality"
$ws.Range("H10").Value = "afmaken"
$ws.Range("A10").Value = "common_pitfalls.ipynb"
$ws.Range("C10").Value = 1
$ws.Range("F10").Value = "Onno Ebbens"

# New row 11: Quality_ZEN_of_Python.ipynb notebook
$ws.Range("A11").Value = "Quality_ZEN_of_Python.ipynb"
$ws.Range("H11").Value = "af"
$ws.Range("E11").Value = "the ZEN of Python"
$ws.Range("B11").Value = "code quality"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = "basis"
$ws.Range("F11").Value = "Onno Ebbens"

# Wrap text formatting to match the rest of the table
$ws.Range("B10:C10").WrapText = $true
$ws.Range("H10").WrapText = $true
$ws.Range("B11:C11").WrapText = $true
$ws.Range("E11").WrapText = $true
$ws.Range("H11").WrapText = $true

# Update active selection to the new bottom of the table
$ws.Range("G11").Select()
